$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 4) below the existing header row (row 3):
# C=ID, D=Id Colaborador, E=Nome Colaborador, F=Mes, G=Horas Mes,
# H=GP, I=Horas Trabalhadas, J=Proporcao de Hora, K=Valor por GP,
# L=GP 9014, M=Observacao GP 9010, N=Observacao GP 9021
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Marlon Passeri"
$ws.Range("F4").Value = "2025-06"
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 656
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 1
# K4, L4, M4, N4 stay empty (only formatted, no value)

# Borders: give the new row the same "bottom of the table" look as the
# header row already has on its sides - a thin bottom border across the
# whole row, plus the outer left/right edges on the first/last columns.
$ws.Range("C4:N4").Borders(9).LineStyle = 1
$ws.Range("C4").Borders(7).LineStyle = 1
$ws.Range("N4").Borders(10).LineStyle = 1

# Number formats for the numeric "hours" / "proportion" columns
$ws.Range("G4").NumberFormat = "0.00"
$ws.Range("I4").NumberFormat = "0.00"
$ws.Range("J4").NumberFormat = "0.00%"
